$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two groups of "sonuc" (result) columns (C/D/E) for the
# temperature rows. Rows 2-5 (temps 5/10/15/20) get the "below 25" results;
# rows 6-8 (temps 25/30/35) get the "25 or above" results. Fill column by
# column (both groups) to match the original authoring order.
$ws.Range("C2:C5").Value = "normal"
$ws.Range("C6:C8").Value = "havasıcak"

$ws.Range("D2:D5").Value = "sıcak ya da soğuk"
$ws.Range("D6:D8").Value = "hava çok güzel ve ılık"

$ws.Range("E2:E5").Value = "havaya dikkat edin"
$ws.Range("E6:E8").Value = "bugün hava 25 derece"

# Fill in the student info block (Numara / Ad Soyad / Bölüm).
$ws.Range("H2").Value = 20215070055
$ws.Range("H3").Value = "Muhammed Ali Harmancı"
$ws.Range("H4").Value = "Yönetim Bilişim Sistemleri"

# Match the saved selection/scroll state from the edited workbook.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("H4:J4").Select()
